# Update PCB BOM: switch SW1/SW2 tactile buttons to smaller 3x2mm Omron buttons.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 holds the SW1, SW2 tactile switch line item.
# Manufacturer: C&K -> Omron Electronics Inc-EMC Div
# Mfg Part #:   PTS810 SJM 250 SMTR LFS -> B3U-1000P
$ws.Range("B5").Value = "B3U-1000P"
$ws.Range("A5").Value = "Omron Electronics Inc-EMC Div"

# Update the active selection left behind by the edit session.
$ws.Range("E9").Select()
